# Updates cryptos list values (Price + Volume(1h), and a few re-ranked
# Coin/Link pairs) to match the latest scrape, per the commit:
# "Updated cryptos list on Wed Oct 30 14:50:19 UTC 2024 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Only the rows whose Coin/Link/Price/Volume actually changed are listed here.
# A $null field means "leave that column alone for this row".
$updates = @(
    @{ Row = 2; B = $null; C = $null; D = '72.311.24'; E = '  +0.67%  ' },
    @{ Row = 3; B = $null; C = $null; D = '2.710.79'; E = '  +3.07%  ' },
    @{ Row = 4; B = $null; C = $null; D = $null; E = '  +0.06%  ' },
    @{ Row = 5; B = $null; C = $null; D = '600.67'; E = '  -1.00%  ' },
    @{ Row = 6; B = $null; C = $null; D = '176.25'; E = '  -1.69%  ' },
    @{ Row = 7; B = $null; C = $null; D = $null; E = '  +0.01%  ' },
    @{ Row = 8; B = $null; C = $null; D = '0.525'; E = '  -0.22%  ' },
    @{ Row = 9; B = $null; C = $null; D = '2.710.26'; E = '  +3.11%  ' },
    @{ Row = 10; B = $null; C = $null; D = $null; E = '  +0.50%  ' },
    @{ Row = 11; B = $null; C = $null; D = $null; E = '  +2.58%  ' },
    @{ Row = 12; B = $null; C = $null; D = '0.355'; E = '  +2.08%  ' },
    @{ Row = 13; B = $null; C = $null; D = '5.02'; E = '  -0.13%  ' },
    @{ Row = 14; B = $null; C = $null; D = '3.206.21'; E = '  +2.29%  ' },
    @{ Row = 15; B = $null; C = $null; D = $null; E = '  -0.02%  ' },
    @{ Row = 16; B = $null; C = $null; D = '72.166.62'; E = '  +0.65%  ' },
    @{ Row = 17; B = $null; C = $null; D = '26.39'; E = '  -0.57%  ' },
    @{ Row = 18; B = $null; C = $null; D = '2.719.03'; E = '  +2.66%  ' },
    @{ Row = 19; B = $null; C = $null; D = '12.31'; E = '  +7.12%  ' },
    @{ Row = 20; B = $null; C = $null; D = '8.15'; E = '  +2.43%  ' },
    @{ Row = 21; B = $null; C = $null; D = '373.81'; E = '  -2.82%  ' },
    @{ Row = 22; B = $null; C = $null; D = $null; E = '  +0.45%  ' },
    @{ Row = 23; B = $null; C = $null; D = $null; E = '  +2.71%  ' },
    @{ Row = 24; B = $null; C = $null; D = '72.49'; E = '  -0.37%  ' },
    @{ Row = 25; B = 'Dai'; C = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'; D = '1.00'; E = '  -0.10%  ' },
    @{ Row = 26; B = 'NEARProtocol'; C = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D = '4.39'; E = '  -1.60%  ' },
    @{ Row = 27; B = $null; C = $null; D = '9.85'; E = '  -0.56%  ' },
    @{ Row = 28; B = $null; C = $null; D = '2.853.71'; E = '  +3.11%  ' },
    @{ Row = 29; B = $null; C = $null; D = '1.00'; E = '  +0.03%  ' },
    @{ Row = 30; B = $null; C = $null; D = '0.0₃0995'; E = '  +3.26%  ' },
    @{ Row = 31; B = $null; C = $null; D = '8.14'; E = '  +1.19%  ' },
    @{ Row = 32; B = $null; C = $null; D = '509.04'; E = '  -6.54%  ' },
    @{ Row = 33; B = $null; C = $null; D = $null; E = '  -1.10%  ' },
    @{ Row = 34; B = $null; C = $null; D = $null; E = '  -0.03%  ' },
    @{ Row = 35; B = $null; C = $null; D = $null; E = '  -0.03%  ' },
    @{ Row = 36; B = $null; C = $null; D = '164.12'; E = '  -1.24%  ' },
    @{ Row = 37; B = $null; C = $null; D = '19.71'; E = '  +2.47%  ' },
    @{ Row = 38; B = $null; C = $null; D = $null; E = '  -0.15%  ' },
    @{ Row = 39; B = $null; C = $null; D = $null; E = '  -0.11%  ' },
    @{ Row = 40; B = $null; C = $null; D = '0.109'; E = '  -4.52%  ' },
    @{ Row = 41; B = $null; C = $null; D = '1.81'; E = '  -2.99%  ' },
    @{ Row = 42; B = $null; C = $null; D = '5.08'; E = '  +0.91%  ' },
    @{ Row = 43; B = $null; C = $null; D = $null; E = '  +0.01%  ' },
    @{ Row = 44; B = 'PolygonEcosystemToken'; C = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'; D = '0.336'; E = '  +1.13%  ' },
    @{ Row = 45; B = 'dogwifhat'; C = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'; D = '2.57'; E = '  -2.16%  ' },
    @{ Row = 46; B = $null; C = $null; D = '156.74'; E = '  +3.75%  ' },
    @{ Row = 47; B = $null; C = $null; D = '39.52'; E = '  +0.74%  ' },
    @{ Row = 48; B = $null; C = $null; D = '0.567'; E = '  +5.89%  ' },
    @{ Row = 49; B = 'Filecoin'; C = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D = '3.76'; E = '  +3.12%  ' },
    @{ Row = 50; B = 'Optimism'; C = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'; D = '1.78'; E = '  +4.98%  ' },
    @{ Row = 51; B = $null; C = $null; D = '0.0766'; E = '  +1.09%  ' }
)

foreach ($u in $updates) {
    $r = $u.Row

    if ($null -ne $u.B) {
        $ws.Range("B$r").Value = $u.B
    }
    if ($null -ne $u.C) {
        $ws.Range("C$r").Value = $u.C
    }
    if ($null -ne $u.D) {
        # Force text formatting first so price strings like "1.00" or
        # "176.25" aren't silently coerced into numbers by Excel.
        $ws.Range("D$r").NumberFormat = "@"
        $ws.Range("D$r").Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Range("E$r").Value = $u.E
    }
}
